# BAU Guaranteed Dispatch Perc by Elec Source.xlsx - apply commit changes
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("BGDPbES")

# --- 1. "About" sheet: remove the old biomass/geothermal/hydro/nuclear notes
#        (rows 13, 14, 16, 17). Delete from the bottom up so row numbers of
#        earlier rows are not disturbed while we work.
$ws1.Rows.Item(17).Delete()
$ws1.Rows.Item(16).Delete()
$ws1.Rows.Item(14).Delete()
$ws1.Rows.Item(13).Delete()

# --- 2. "BGDPbES" sheet: fix rows 9 (biomass) & 10 (geothermal), which had
#        hard-coded 1's instead of following the sheet's "=$B<row>" fill
#        pattern used by every other row. Set B to 0 and re-fill C:AK.
$ws2.Range("B9").Value = 0
$ws2.Range("C9:AK9").Formula = '=$B9'

$ws2.Range("B10").Value = 0
$ws2.Range("C10:AK10").Formula = '=$B10'

# --- 3. Add three new fuel-source rows (crude oil, heavy/residual fuel oil,
#        municipal solid waste), all with zero values across every year.
#        (Done before the A1 header text below so the new shared-string
#        table entries line up in the same order as the source workbook.)
$ws2.Range("A15").Value = "crude oil"
$ws2.Range("B15:AK15").Value = 0

$ws2.Range("A16").Value = "heavy or residual fuel oil"
$ws2.Range("B16:AK16").Value = 0

$ws2.Range("A17").Value = "municipal solid waste"
$ws2.Range("B17:AK17").Value = 0

# --- 4. Add a bold, wrapped header label in A1 describing the table, and
#        give row 1 extra height so the wrapped text is visible.
$ws2.Range("A1").Value = "BAU Guaranteed Dispatch (dimensionless)"
$ws2.Range("A1").Font.Bold = $true
$ws2.Range("A1").WrapText = $true
$ws2.Rows.Item(1).RowHeight = 45

# --- 5. Misc page setup matching the new layout.
$ws2.PageSetup.Orientation = 1

Write-Host "done"
